# Update betting odds values in row 4 of the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G4"  = 5.8
    "H4"  = 4.5
    "I4"  = 1.47
    "J4"  = 5.4
    "K4"  = 2.47
    "L4"  = 1.93
    "U4"  = 1.7
    "V4"  = 2.02
    "W4"  = 17
    "X4"  = 40
    "Y4"  = 19.5
    "Z4"  = 120
    "AA4" = 60
    "AB4" = 55
    "AD4" = 9.5
    "AE4" = 17.5
    "AF4" = 70
    "AG4" = 500
    "AH4" = 8
    "AI4" = 8.5
    "AK4" = 11.5
    "AL4" = 11.75
    "AM4" = 23
    "AN4" = 7.4
    "AO4" = 30
    "AP4" = 30
    "AQ4" = 175
    "AR4" = 175
    "AS4" = 350
    "AU4" = 7.5
    "AV4" = 60
    "AW4" = 3.5
    "AX4" = 6.7
    "AZ4" = 18.5
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
